$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 221 is being bulk-updated
# from date serial 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C221").Value = 45175
